# Updated the valid email address
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC02_UserProfile")

# --- Update the data row (row 2) values ---
# (order matters for shared-string table layout: test, then auto, then password@1)
$ws.Range("B2").Value = "test"
$ws.Range("A2").Value = "auto"
$ws.Range("D2").Value = "password@1"
$ws.Range("E2").Value = "password@1"
$ws.Range("F2").Value = 11111

# --- Rebuild the hyperlinks so that the resulting order is D2, K2, C2, E2 ---
# (the underlying mailto targets themselves are unchanged)
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:dada@123")
$ws.Range("D2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:dada@123")
$ws.Range("K2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:dada@1234", "", "", "venkat@cgi.com")
$ws.Range("C2").Style = "Hyperlink"
# Adding the hyperlink overwrote C2's formula with its display text - restore it.
$ws.Range("C2").Formula = "=CONCATENATE(A2,L2)"

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:venkat@cgi.com")
$ws.Range("E2").Style = "Hyperlink"

# --- Update the sheet view: drop the scrolled topLeftCell and move the selection to A2 ---
$ws.Range("A2").Select()
